$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 121
$ws.Range("B3").Value = 104

$ws.Range("A4").Value = "GUTIERREZ CARLOS TERESA DE JESUS"
$ws.Range("B4").Value = 103

$ws.Range("A5").Value = "ZAVALETA MANAY JORGE LUIS"
$ws.Range("B5").Value = 96

$ws.Range("B6").Value = 95

$ws.Range("A7").Value = "HIDALGO CUBAS LUISA YVONE"
$ws.Range("B7").Value = 95

$ws.Range("A8").Value = "SEVERINO AVALOS MARJORIE ISABEL"
$ws.Range("B8").Value = 91

$ws.Range("A9").Value = "SENADOR ARBOLEDA GIANCARLOS EXEBIO"
$ws.Range("B9").Value = 88

$ws.Range("A10").Value = "HUMPIRE CASTILLO IRWIN DEIMER"
$ws.Range("B10").Value = 87

$ws.Range("A11").Value = "BALLENA ESQUÉN ASTRID CAROLINA"
$ws.Range("B11").Value = 87

$ws.Range("A12").Value = "VALLE MAGALLAN EDUAR"
$ws.Range("B12").Value = 86

$ws.Range("A13").Value = "CAMACHO LINARES JUDITH ARLETT"
$ws.Range("B13").Value = 76
